$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "2025/09/30"
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").Value = "火"
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 16
